# Append two new time-log entries to the sheet (rows 52 and 53),
# matching the style/number-format already used for the date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date number format from the last populated row (A51)
# so the new date cells get the same cell style ("s") instead of a brand
# new custom numFmt.
$dateFmt = $ws.Cells.Item(51, 1).NumberFormat

# Row 52: 11/14/2023, 2 hours, debugging note
$ws.Cells.Item(52, 1).Value = 45244
$ws.Cells.Item(52, 1).NumberFormat = $dateFmt
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(52, 3).Value = "bebugging ie some of the FXML and the controller integration class werent working properly, and data wasn’t being printed"

# Row 53: 11/16/2023, 3 hours, frontend/backend integration note
$ws.Cells.Item(53, 1).Value = 45246
$ws.Cells.Item(53, 1).NumberFormat = $dateFmt
$ws.Cells.Item(53, 2).Value = 3
$ws.Cells.Item(53, 3).Value = "starting the back-end and frontend UI integrations with the APIS, and front end"

# Match the saved selection state from the workbook (active cell C53)
[void]$ws.Range("C53").Select()
